$d = $word.ActiveDocument

function Find-RangeByText($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "text not found: $searchText"
    }
    return $rng
}

# ---------------------------------------------------------------------------
# 1) "Pouldre d'horloges de sable" -> "Pouldre d" + "'" + "horloges de sable"
#    the apostrophe becomes a curly quote in its own (non-colored) run
# ---------------------------------------------------------------------------
$whole = Find-RangeByText "Pouldre d'horloges de sable"
$base = $whole.Start
$aposStart = $base + 9
$aposEnd = $aposStart + 1
$rApos = $d.Range($aposStart, $aposEnd)
$rApos.Text = [char]0x2019
$rApos.Font.Color = -16777216

# ---------------------------------------------------------------------------
# 2) "purifier de sa crasse, puys versé dedans quattre " ->
#    "purifier de sa crasse, puys vers" + "e" + " dedans quattre "
#    the accented e becomes a plain e in its own (non-colored) run
# ---------------------------------------------------------------------------
$whole = Find-RangeByText "purifier de sa crasse, puys versé dedans quattre "
$base = $whole.Start
$eStart = $base + 32
$eEnd = $eStart + 1
$rE = $d.Range($eStart, $eEnd)
$rE.Text = "e"
$rE.Font.Color = -16777216

# ---------------------------------------------------------------------------
# 3) "incorporé, et leve le incontinent du feu tousjours meslant. Et s'il te"
#    -> same text but with a comma inserted after "feu" (single run, no split)
# ---------------------------------------------------------------------------
$whole = Find-RangeByText "incorporé, et leve le incontinent du feu tousjours meslant. Et s'il te"
$whole.Text = "incorporé, et leve le incontinent du feu, tousjours meslant. Et s'il te"

# ---------------------------------------------------------------------------
# 4) " fin," -> " fin" + "."   (comma replaced by period, in its own run)
# ---------------------------------------------------------------------------
$whole = Find-RangeByText " fin,"
$base = $whole.Start
$commaStart = $base + 4
$commaEnd = $commaStart + 1
$rComma = $d.Range($commaStart, $commaEnd)
$rComma.Text = "."
$rComma.Font.Color = -16777216

# ---------------------------------------------------------------------------
# 5) "puys lave le tant de fois que l'" -> "P" + "uys lave le tant de fois que l'"
#    (capitalize the first letter, in its own non-colored run)
# ---------------------------------------------------------------------------
$whole = Find-RangeByText "puys lave le tant de fois que l'"
$base = $whole.Start
$pStart = $base
$pEnd = $pStart + 1
$rP = $d.Range($pStart, $pEnd)
$rP.Text = "P"
$rP.Font.Color = -16777216

# ---------------------------------------------------------------------------
# 6) "t ceste petite" -> "t ceste " + "petite"  (split, same formatting kept)
# ---------------------------------------------------------------------------
$whole = Find-RangeByText "t ceste petite"
$base = $whole.Start
$splitPos = $base + 8
$rFirst = $d.Range($base, $splitPos)
$origColor = $rFirst.Font.Color
$rFirst.Font.Color = $origColor + 1
$rFirst.Font.Color = $origColor
